$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is numeric-looking text (e.g. "1.000", "0.7703")
# must be forced to the Text number format first, otherwise Excel COM
# auto-coerces the assignment into a real number (losing e.g. trailing
# zeros) -- the source data stores every value as text.
$ws.Range("D2").Value = '29.872.46'
$ws.Range("E2").Value = '  -0.08%  '
$ws.Range("D3").Value = '1.888.14'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7703'
$ws.Range("E5").Value = '  -0.84%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.54'
$ws.Range("E6").Value = '  -0.98%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3119'
$ws.Range("E8").Value = '  -0.99%  '
$ws.Range("E9").Value = '  +0.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07179'
$ws.Range("E10").Value = '  -6.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08601'
$ws.Range("E11").Value = '  +6.08%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7639'
$ws.Range("E12").Value = '  -1.20%  '
$ws.Range("D13").Value = '1.902.04'
$ws.Range("E13").Value = '  +2.77%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.363'
$ws.Range("E14").Value = '  -2.48%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '93.57'
$ws.Range("E15").Value = '  +1.13%  '
$ws.Range("E16").Value = '  -1.65%  '
$ws.Range("D17").Value = '29.870.81'
$ws.Range("E17").Value = '  -0.03%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.76'
$ws.Range("E18").Value = '  -1.93%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.55'
$ws.Range("E19").Value = '  -0.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007806'
$ws.Range("E20").Value = '  -1.98%  '
$ws.Range("D21").Value = '2.155.67'
$ws.Range("E21").Value = '  +1.84%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9998'
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.025'
$ws.Range("E23").Value = '  -1.63%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1646'
$ws.Range("E25").Value = '  +4.61%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.373'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.61'
$ws.Range("E27").Value = '  -0.21%  '
$ws.Range("E28").Value = '  -0.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.034'
$ws.Range("E29").Value = '  -0.77%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.461'
$ws.Range("E30").Value = '  +1.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.536'
$ws.Range("E31").Value = '  -1.03%  '
$ws.Range("E32").Value = '  +0.08%  '
$ws.Range("E33").Value = '  -0.18%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05446'
$ws.Range("E34").Value = '  -1.55%  '
$ws.Range("E35").Value = '  -1.85%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7440'
$ws.Range("E36").Value = '  -2.09%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.001'
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.701'
$ws.Range("E38").Value = '  +2.23%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01954'
$ws.Range("E39").Value = '  +1.19%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.782'
$ws.Range("E40").Value = '  -0.23%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4466'
$ws.Range("E41").Value = '  +0.44%  '
$ws.Range("D42").Value = '1.109.46'
$ws.Range("E42").Value = '  -4.26%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.072'
$ws.Range("E43").Value = '  +2.09%  '
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '72.97'
$ws.Range("E44").Value = '  -1.55%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8511'
$ws.Range("E45").Value = '  +0.38%  '
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.38'
$ws.Range("E47").Value = '  -0.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.862'
$ws.Range("E48").Value = '  -2.05%  '
$ws.Range("E49").Value = '  +1.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.014'
$ws.Range("E50").Value = '  -4.16%  '
$ws.Range("D51").Value = '2.056.42'
$ws.Range("E51").Value = '  +1.05%  '
